$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update odds values for rows 2-7 (minor odds adjustments)
$ws.Range("M2").Value = 1.03
$ws.Range("O2").Value = 1.18

$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.36
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83

$ws.Range("M4").Value = 1.06
$ws.Range("O4").Value = 1.3
$ws.Range("U4").Value = 1.73

$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.29
$ws.Range("U5").Value = 1.67

$ws.Range("M6").Value = 1.04
$ws.Range("O6").Value = 1.25
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.1
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91

$ws.Range("M7").Value = 1.05
$ws.Range("O7").Value = 1.29
$ws.Range("U7").Value = 1.73

# 2) Remove row 8 (ESTONIA - MEISTRILIIGA, Flora vs Kalju) entirely,
#    shifting all subsequent rows up by one.
$ws.Rows.Item(8).Delete()
